# ApachePOI cucumber scenario - steps -2
# 1. Rename the shared strings used on the "testCitizen" sheet
#    (columns A and B, rows 1-8) to their new values.
# 2. Move the active-tab / tab-selected state from Sheet1 to testCitizen
#    and update the selected cell on testCitizen.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)
$sheet2 = $wb.Worksheets.Item(2)

$newColA = @("ulis111","ulis112","ulis113","ulis114","ulis115","ulis116","ulis117","ulis118")
$newColB = @("ubs111","ubs112","ubs113","ubs114","ubs115","ubs116","ubs117","ubs118")

for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 1
    $sheet2.Cells.Item($row, 1).Value = $newColA[$i]
    $sheet2.Cells.Item($row, 2).Value = $newColB[$i]
}

# Switch the active sheet/tab to testCitizen and select F14 there.
$sheet2.Activate()
$sheet2.Range("F14").Select()

# Keep Sheet1's own remembered selection as it was (A5).
$sheet1.Range("A5").Select()

# Re-activate testCitizen so it ends up as the workbook's active tab.
$sheet2.Activate()
